# Auto-generated edit script applying numeric corrections to LeveProfit-related
# columns (H..N) across the per-job worksheets, per scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 86.44444
$ws.Range("I4").Value = 86.44444
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 86.44444
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 27.55556
$ws.Range("N4").ClearContents()
$ws.Range("H12").Value = 750.25
$ws.Range("I12").Value = 550.5
$ws.Range("K12").Value = 550.5
$ws.Range("M12").Value = -380.5
$ws.Range("H18").Value = 1716.3334
$ws.Range("I18").Value = 1099.5385
$ws.Range("J18").Value = 3320
$ws.Range("K18").Value = 1099.5385
$ws.Range("L18").Value = 3320
$ws.Range("M18").Value = -815.5385000000001
$ws.Range("N18").Value = -3888
$ws.Range("H43").Value = 8944.333000000001
$ws.Range("I43").Value = 9000
$ws.Range("K43").Value = 9000
$ws.Range("M43").Value = -8931
$ws.Range("H53").Value = 597.6111
$ws.Range("I53").Value = 392.33334
$ws.Range("J53").Value = 700.25
$ws.Range("K53").Value = 392.33334
$ws.Range("L53").Value = 700.25
$ws.Range("M53").Value = 244.66666
$ws.Range("N53").Value = -1974.25
$ws.Range("H86").Value = 22325456
$ws.Range("I86").Value = 31251000
$ws.Range("J86").Value = 20837866
$ws.Range("K86").Value = 31251000
$ws.Range("L86").Value = 20837866
$ws.Range("M86").Value = -31249877
$ws.Range("N86").Value = -20840112
$ws.Range("H89").Value = 22325456
$ws.Range("I89").Value = 31251000
$ws.Range("J89").Value = 20837866
$ws.Range("K89").Value = 156255000
$ws.Range("L89").Value = 104189330
$ws.Range("M89").Value = -156249384
$ws.Range("N89").Value = -104200562
$ws.Range("H92").Value = 950.2105
$ws.Range("I92").Value = 697.1875
$ws.Range("K92").Value = 697.1875
$ws.Range("M92").Value = 550.8125
$ws.Range("H106").Value = 4424.75
$ws.Range("I106").Value = 3598.8
$ws.Range("K106").Value = 3598.8
$ws.Range("M106").Value = -2967.8
$ws.Range("H112").Value = 1856.5405
$ws.Range("J112").Value = 1929.2059
$ws.Range("L112").Value = 5787.6177
$ws.Range("N112").Value = -8003.6177
$ws.Range("H137").Value = 3599.879
$ws.Range("I137").Value = 1710.421
$ws.Range("K137").Value = 5131.263
$ws.Range("M137").Value = -2581.263

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1329.7966
$ws.Range("I32").Value = 942.7455
$ws.Range("K32").Value = 942.7455
$ws.Range("M32").Value = -655.7455
$ws.Range("H74").Value = 7344.2
$ws.Range("I74").Value = 2246
$ws.Range("K74").Value = 2246
$ws.Range("M74").Value = -1372
$ws.Range("H77").Value = 7344.2
$ws.Range("I77").Value = 2246
$ws.Range("K77").Value = 11230
$ws.Range("M77").Value = -6862
$ws.Range("H132").Value = 2350.875
$ws.Range("I132").Value = 2168.6453
$ws.Range("K132").Value = 6505.9359
$ws.Range("M132").Value = -3975.9359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 527.53845
$ws.Range("I94").Value = 471.5
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 471.5
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -20.5
$ws.Range("N94").Value = -2102
$ws.Range("H107").Value = 4143.0713
$ws.Range("I107").Value = 4334.4443
$ws.Range("J107").Value = 3798.6
$ws.Range("K107").Value = 4334.4443
$ws.Range("L107").Value = 3798.6
$ws.Range("M107").Value = -2414.4443
$ws.Range("N107").Value = -7638.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1290.4
$ws.Range("I5").Value = 1489.25
$ws.Range("J5").Value = 495
$ws.Range("K5").Value = 1489.25
$ws.Range("L5").Value = 495
$ws.Range("M5").Value = -1377.25
$ws.Range("N5").Value = -719
$ws.Range("H10").Value = 271.33334
$ws.Range("I10").Value = 271.33334
$ws.Range("K10").Value = 271.33334
$ws.Range("M10").Value = -132.33334
$ws.Range("H25").Value = 10402.2
$ws.Range("I25").Value = 502.75
$ws.Range("J25").Value = 50000
$ws.Range("K25").Value = 502.75
$ws.Range("L25").Value = 50000
$ws.Range("M25").Value = -328.75
$ws.Range("N25").Value = -50348
$ws.Range("H31").Value = 8677.540999999999
$ws.Range("I31").Value = 1068.1428
$ws.Range("K31").Value = 1068.1428
$ws.Range("M31").Value = -773.1428000000001
$ws.Range("H34").Value = 8677.540999999999
$ws.Range("I34").Value = 1068.1428
$ws.Range("K34").Value = 1068.1428
$ws.Range("M34").Value = -866.1428000000001
$ws.Range("H58").Value = 12433.8
$ws.Range("I58").Value = 24131.182
$ws.Range("K58").Value = 24131.182
$ws.Range("M58").Value = -23928.182
$ws.Range("H136").Value = 12433.8
$ws.Range("I136").Value = 24131.182
$ws.Range("K136").Value = 72393.546
$ws.Range("M136").Value = -69843.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 54883120
$ws.Range("I4").Value = 85542750
$ws.Range("K4").Value = 256628250
$ws.Range("M4").Value = -256628138
$ws.Range("H38").Value = 22
$ws.Range("I38").Value = 4
$ws.Range("K38").Value = 12
$ws.Range("M38").Value = 335
$ws.Range("H68").Value = 11755.212
$ws.Range("I68").Value = 20618.666
$ws.Range("J68").Value = 6690.381
$ws.Range("K68").Value = 61855.99800000001
$ws.Range("L68").Value = 20071.143
$ws.Range("M68").Value = -61044.99800000001
$ws.Range("N68").Value = -21693.143
$ws.Range("H71").Value = 11755.212
$ws.Range("I71").Value = 20618.666
$ws.Range("J71").Value = 6690.381
$ws.Range("K71").Value = 185567.994
$ws.Range("L71").Value = 60213.429
$ws.Range("M71").Value = -181511.994
$ws.Range("N71").Value = -68325.429
$ws.Range("H122").Value = 12500366
$ws.Range("I122").Value = 377.25
$ws.Range("J122").Value = 50000332
$ws.Range("K122").Value = 3395.25
$ws.Range("L122").Value = 450002988
$ws.Range("M122").Value = -945.25
$ws.Range("N122").Value = -450007888
$ws.Range("H131").Value = 2164.96
$ws.Range("I131").Value = 832.8461
$ws.Range("J131").Value = 3608.0833
$ws.Range("K131").Value = 2498.5383
$ws.Range("L131").Value = 10824.2499
$ws.Range("M131").Value = 2541.4617
$ws.Range("N131").Value = -20904.2499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 15000
$ws.Range("I35").Value = 15000
$ws.Range("K35").Value = 15000
$ws.Range("M35").Value = -14702
$ws.Range("H43").Value = 5939.875
$ws.Range("J43").Value = 17009.5
$ws.Range("L43").Value = 17009.5
$ws.Range("N43").Value = -17311.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3365.6584
$ws.Range("I16").Value = 3857.6572
$ws.Range("J16").Value = 495.66666
$ws.Range("K16").Value = 3857.6572
$ws.Range("L16").Value = 495.66666
$ws.Range("M16").Value = -3687.6572
$ws.Range("N16").Value = -835.66666
$ws.Range("H64").Value = 16798
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 19747.5
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 19747.5
$ws.Range("M64").Value = -4775
$ws.Range("N64").Value = -20197.5
$ws.Range("H67").Value = 16798
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 19747.5
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 19747.5
$ws.Range("M67").Value = -4220
$ws.Range("N67").Value = -21307.5
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 1000
$ws.Range("M93").Value = 248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1401.3334
$ws.Range("I7").Value = 2002
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 2002
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -1889
$ws.Range("N7").Value = -426
$ws.Range("H68").Value = 80012
$ws.Range("J68").Value = 80012
$ws.Range("L68").Value = 80012
$ws.Range("N68").Value = -81634
$ws.Range("H71").Value = 80012
$ws.Range("J71").Value = 80012
$ws.Range("L71").Value = 240036
$ws.Range("N71").Value = -248148
$ws.Range("H86").Value = 51464.4
$ws.Range("J86").Value = 51464.4
$ws.Range("L86").Value = 51464.4
$ws.Range("N86").Value = -53710.4
$ws.Range("H89").Value = 51464.4
$ws.Range("J89").Value = 51464.4
$ws.Range("L89").Value = 257322
$ws.Range("N89").Value = -268554
$ws.Range("I107").Value = 1124.75
$ws.Range("J107").Value = 4497.5
$ws.Range("K107").Value = 3374.25
$ws.Range("L107").Value = 13492.5
$ws.Range("M107").Value = -1454.25
$ws.Range("N107").Value = -17332.5
$ws.Range("H132").Value = 6599.2354
$ws.Range("I132").Value = 2932.75
$ws.Range("K132").Value = 8798.25
$ws.Range("M132").Value = -6268.25
